# Eukardia Πεδία Βάσης Δεδομένων - add "SQL view" calculated column to Πίνακας3 (sheet "dim")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dim")
$ws.Activate()

$lo = $ws.ListObjects.Item(1)

# Add a new (8th) column to the table - this extends the table ref from A1:G42 to A1:H42
$newCol = $lo.ListColumns.Add()

# Header text for the new column
$ws.Range("H1").Value = "SQL view"

# Same calculated-column formula used for every data row, written cell-by-cell
# (matches the workbook's existing convention of per-cell literal formulas rather
# than a single shared formula block).
$formula = '="`cases`.`"&Πίνακας3[[#This Row],[Πεδίο]]&"` AS `"&Πίνακας3[[#This Row],[Πεδίο]]&"`,"'

for ($r = 2; $r -le 42; $r++) {
    $cell = $ws.Range("H$r")
    $cell.Formula = $formula
}

# Match the style used by the rest of the table's data rows (vertical-center + wrap)
$dataRng = $ws.Range("H2:H42")
$dataRng.WrapText = $true
$dataRng.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# Widen the new column like the others (bestFit-style width)
$ws.Columns.Item(8).ColumnWidth = 40.45

# Reflect the author's selection / scroll position when they added the column
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("H7:H41").Select()
